# Updates cryptos list values (Price and Volume(1h) columns, and a few Coin/Link
# cells where rankings shifted) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '71.147.98'
$ws.Range('E2').Value = '  +2.58%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.688.80'
$ws.Range('E3').Value = '  +7.91%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.22'
$ws.Range('E5').Value = '  -0.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.45'
$ws.Range('E6').Value = '  +0.66%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.676.99'
$ws.Range('E7').Value = '  +7.78%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.615'
$ws.Range('E8').Value = '  +4.12%  '
$ws.Range('E9').Value = '  +0.11%  '
$ws.Range('E10').Value = '  +0.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.74'
$ws.Range('E11').Value = '  +24.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.608'
$ws.Range('E12').Value = '  +4.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '48.87'
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('E14').Value = '  +2.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.279.98'
$ws.Range('E15').Value = '  +7.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '677.59'
$ws.Range('E16').Value = '  -1.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '8.99'
$ws.Range('E17').Value = '  +4.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.680.26'
$ws.Range('E18').Value = '  +7.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.210.81'
$ws.Range('E19').Value = '  +2.56%  '
$ws.Range('E20').Value = '  +1.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.91'
$ws.Range('E21').Value = '  +1.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.57'
$ws.Range('E22').Value = '  +2.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.939'
$ws.Range('E23').Value = '  +5.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '17.32'
$ws.Range('E24').Value = '  +2.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '101.88'
$ws.Range('E25').Value = '  +0.91%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.96'
$ws.Range('E26').Value = '  +1.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.82'
$ws.Range('E27').Value = '  +5.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.25'
$ws.Range('E28').Value = '  +7.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '35.04'
$ws.Range('E29').Value = '  +5.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.17'
$ws.Range('E30').Value = '  +5.35%  '
$ws.Range('B31').Value = 'Mantle'
$ws.Range('C31').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.42'
$ws.Range('E31').Value = '  -1.88%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.47'
$ws.Range('E32').Value = '  +6.29%  '
$ws.Range('B33').Value = 'dogwifhat'
$ws.Range('C33').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.08'
$ws.Range('E33').Value = '  +11.75%  '
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '586.19'
$ws.Range('E34').Value = '  +1.71%  '
$ws.Range('B35').Value = 'Cosmos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '11.17'
$ws.Range('E35').Value = '  +1.73%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.108'
$ws.Range('E36').Value = '  +5.35%  '
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '58.93'
$ws.Range('E37').Value = '  +1.13%  '
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.680.45'
$ws.Range('E39').Value = '  +3.12%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.145'
$ws.Range('E41').Value = '  +4.51%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '35.32'
$ws.Range('E42').Value = '  +1.59%  '
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0₃0764'
$ws.Range('E43').Value = '  +5.12%  '
$ws.Range('B44').Value = 'Stacks'
$ws.Range('C44').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.41'
$ws.Range('E44').Value = '  +5.45%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.75'
$ws.Range('E45').Value = '  +3.85%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0453'
$ws.Range('E46').Value = '  +9.10%  '
$ws.Range('B47').Value = 'TheGraph'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.348'
$ws.Range('E47').Value = '  +5.18%  '
$ws.Range('B48').Value = 'ThetaToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.88'
$ws.Range('E48').Value = '  +9.04%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.34'
$ws.Range('E49').Value = '  -0.61%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.133'
$ws.Range('E50').Value = '  +3.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '135.97'
$ws.Range('E51').Value = '  +2.71%  '

Write-Host "Applied cryptos update"
